$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 13274
$ws.Range("E2").Value = 816
$ws.Range("F2").Value = 816
$ws.Range("G2").Value = 5796
$ws.Range("H2").Value = 5576
$ws.Range("I2").Value = 5490
$ws.Range("J2").Value = 86
$ws.Range("K2").Value = 355074
$ws.Range("L2").Value = 332280
$ws.Range("M2").Value = 22794
$ws.Range("N2").Value = 18684
$ws.Range("O2").Value = 4111
$ws.Range("P2").Value = 6325
$ws.Range("Q2").Value = -11985
$ws.Range("R2").Value = -464
$ws.Range("S2").Value = 16381
$ws.Range("T2").Value = 272
$ws.Range("U2").Value = $null
$ws.Range("V2").Value = 83952
$ws.Range("W2").Value = 6.15
$ws.Range("X2").Value = 42
$ws.Range("Y2").Value = 41.34
$ws.Range("Z2").Value = 2.16
$ws.Range("AA2").Value = 1457.73
$ws.Range("AB2").Value = 261.36
$ws.Range("AC2").Value = 5259
$ws.Range("AD2").Value = 1.16
$ws.Range("AE2").Value = 14769
$ws.Range("AF2").Value = 0.41
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 1.64
$ws.Range("AI2").Value = 2.3
$ws.Range("AJ2").Value = 126503947

# Row 3
$ws.Range("D3").Value = 19949
$ws.Range("E3").Value = 1964
$ws.Range("F3").Value = 1964
$ws.Range("G3").Value = 1993
$ws.Range("H3").Value = 1509
$ws.Range("I3").Value = 1147
$ws.Range("J3").Value = 362
$ws.Range("K3").Value = 398112
$ws.Range("L3").Value = 371324
$ws.Range("M3").Value = 26787
$ws.Range("N3").Value = 21399
$ws.Range("O3").Value = 5389
$ws.Range("P3").Value = 7772
$ws.Range("Q3").Value = -7449
$ws.Range("R3").Value = -570
$ws.Range("S3").Value = 6560
$ws.Range("T3").Value = 550
$ws.Range("U3").Value = $null
$ws.Range("V3").Value = 88702
$ws.Range("W3").Value = 9.85
$ws.Range("X3").Value = 7.56
$ws.Range("Y3").Value = 5.72
$ws.Range("Z3").Value = 0.4
$ws.Range("AA3").Value = 1386.19
$ws.Range("AB3").Value = 245.47
$ws.Range("AC3").Value = 881
$ws.Range("AD3").Value = 6.24
$ws.Range("AE3").Value = 13767
$ws.Range("AF3").Value = 0.4
$ws.Range("AG3").Value = 50
$ws.Range("AH3").Value = 0.91
$ws.Range("AI3").Value = 6.78
$ws.Range("AJ3").Value = 155439423

# Row 4
$ws.Range("D4").Value = 20971
$ws.Range("E4").Value = 2527
$ws.Range("F4").Value = 2527
$ws.Range("G4").Value = 2620
$ws.Range("H4").Value = 2018
$ws.Range("I4").Value = 1426
$ws.Range("J4").Value = 592
$ws.Range("K4").Value = 457989
$ws.Range("L4").Value = 428351
$ws.Range("M4").Value = 29638
$ws.Range("N4").Value = 22792
$ws.Range("O4").Value = 6846
$ws.Range("P4").Value = 7772
$ws.Range("Q4").Value = -6692
$ws.Range("R4").Value = 21
$ws.Range("S4").Value = 8557
$ws.Range("T4").Value = 395
$ws.Range("U4").Value = $null
$ws.Range("V4").Value = 97284
$ws.Range("W4").Value = 12.05
$ws.Range("X4").Value = 9.63
$ws.Range("Y4").Value = 6.46
$ws.Range("Z4").Value = 0.47
$ws.Range("AA4").Value = 1445.27
$ws.Range("AB4").Value = 282.15
$ws.Range("AC4").Value = 918
$ws.Range("AD4").Value = 6.3
$ws.Range("AE4").Value = 14663
$ws.Range("AF4").Value = 0.39
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 0.87
$ws.Range("AI4").Value = 5.45
$ws.Range("AJ4").Value = 155439423

# Row 5
$ws.Range("D5").Value = 22919
$ws.Range("E5").Value = 3480
$ws.Range("F5").Value = 3480
$ws.Range("G5").Value = 3452
$ws.Range("H5").Value = 2644
$ws.Range("I5").Value = 1851
$ws.Range("J5").Value = 794
$ws.Range("K5").Value = 475937
$ws.Range("L5").Value = 445073
$ws.Range("M5").Value = 30864
$ws.Range("N5").Value = 24295
$ws.Range("O5").Value = 6569
$ws.Range("P5").Value = 7772
$ws.Range("Q5").Value = 11090
$ws.Range("R5").Value = -645
$ws.Range("S5").Value = -12460
$ws.Range("T5").Value = 417
$ws.Range("U5").Value = $null
$ws.Range("V5").Value = 86394
$ws.Range("W5").Value = 15.18
$ws.Range("X5").Value = 11.54
$ws.Range("Y5").Value = 7.86
$ws.Range("Z5").Value = 0.57
$ws.Range("AA5").Value = 1442.05
$ws.Range("AB5").Value = 299.05
$ws.Range("AC5").Value = 1191
$ws.Range("AD5").Value = 5.1
$ws.Range("AE5").Value = 15895
$ws.Range("AF5").Value = 0.38
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 1.65
$ws.Range("AI5").Value = 8.26
$ws.Range("AJ5").Value = 155439423

# Row 6
$ws.Range("D6").Value = 24489
$ws.Range("E6").Value = 4168
$ws.Range("F6").Value = 4168
$ws.Range("G6").Value = 4165
$ws.Range("H6").Value = 3210
$ws.Range("I6").Value = 2415
$ws.Range("K6").Value = 467798
$ws.Range("L6").Value = 433215
$ws.Range("M6").Value = 34583
$ws.Range("N6").Value = 31317
$ws.Range("P6").Value = 9849
$ws.Range("Q6").Value = 4321
$ws.Range("R6").Value = 790
$ws.Range("S6").Value = -6253
$ws.Range("T6").Value = 515
$ws.Range("U6").Value = $null
$ws.Range("V6").Value = 79130
$ws.Range("W6").Value = 17.02
$ws.Range("X6").Value = 13.11
$ws.Range("Y6").Value = 8.69
$ws.Range("Z6").Value = 0.68
$ws.Range("AA6").Value = 1252.68
$ws.Range("AB6").Value = 252.68
$ws.Range("AC6").Value = 1465
$ws.Range("AD6").Value = 3.89
$ws.Range("AE6").Value = 16114
$ws.Range("AF6").Value = 0.35
$ws.Range("AG6").Value = 180
$ws.Range("AH6").Value = 3.16
$ws.Range("AI6").Value = 14.48
$ws.Range("AJ6").Value = 196982894

# Row 7
$ws.Range("D7").Value = $null
$ws.Range("E7").Value = 4754
$ws.Range("G7").Value = 4666
$ws.Range("H7").Value = 3496
$ws.Range("I7").Value = 3343
$ws.Range("K7").Value = 470168
$ws.Range("L7").Value = 432220
$ws.Range("M7").Value = 37948
$ws.Range("N7").Value = 35028
$ws.Range("P7").Value = 9850
$ws.Range("Q7").Value = $null
$ws.Range("R7").Value = $null
$ws.Range("S7").Value = $null
$ws.Range("T7").Value = $null
$ws.Range("U7").Value = $null
$ws.Range("W7").Value = $null
$ws.Range("X7").Value = $null
$ws.Range("Y7").Value = 10.08
$ws.Range("Z7").Value = 0.75
$ws.Range("AA7").Value = 1138.98
$ws.Range("AC7").Value = 1697
$ws.Range("AD7").Value = 3.01
$ws.Range("AE7").Value = 18024
$ws.Range("AF7").Value = 0.28
$ws.Range("AG7").Value = 277
$ws.Range("AH7").Value = 5.41
$ws.Range("AI7").Value = 16.3

# Row 8
$ws.Range("D8").Value = $null
$ws.Range("E8").Value = 4664
$ws.Range("G8").Value = 4639
$ws.Range("H8").Value = 3467
$ws.Range("I8").Value = 3323
$ws.Range("K8").Value = 486818
$ws.Range("L8").Value = 445887
$ws.Range("M8").Value = 40944
$ws.Range("N8").Value = 37772
$ws.Range("P8").Value = 9850
$ws.Range("Q8").Value = $null
$ws.Range("R8").Value = $null
$ws.Range("S8").Value = $null
$ws.Range("T8").Value = $null
$ws.Range("U8").Value = $null
$ws.Range("W8").Value = $null
$ws.Range("X8").Value = $null
$ws.Range("Y8").Value = 9.13
$ws.Range("Z8").Value = 0.73
$ws.Range("AA8").Value = 1089.02
$ws.Range("AC8").Value = 1687
$ws.Range("AD8").Value = 3.03
$ws.Range("AE8").Value = 19435
$ws.Range("AF8").Value = 0.26
$ws.Range("AG8").Value = 306
$ws.Range("AH8").Value = 5.99
$ws.Range("AI8").Value = 18.14

# Row 9
$ws.Range("D9").Value = $null
$ws.Range("E9").Value = 4773
$ws.Range("G9").Value = 4748
$ws.Range("H9").Value = 3543
$ws.Range("I9").Value = 3421
$ws.Range("K9").Value = 503223
$ws.Range("L9").Value = 459316
$ws.Range("M9").Value = 43910
$ws.Range("N9").Value = 40881
$ws.Range("P9").Value = 9850
$ws.Range("Q9").Value = $null
$ws.Range("R9").Value = $null
$ws.Range("S9").Value = $null
$ws.Range("T9").Value = $null
$ws.Range("U9").Value = $null
$ws.Range("W9").Value = $null
$ws.Range("X9").Value = $null
$ws.Range("Y9").Value = 8.7
$ws.Range("Z9").Value = 0.72
$ws.Range("AA9").Value = 1046.03
$ws.Range("AC9").Value = 1737
$ws.Range("AD9").Value = 2.94
$ws.Range("AE9").Value = 21036
$ws.Range("AF9").Value = 0.24
$ws.Range("AG9").Value = 336
$ws.Range("AH9").Value = 6.57
$ws.Range("AI9").Value = 19.33
